# Generate Report for Handback
# The handback transform failed for the db5902a2-... file (both zh-cn and de-de
# targets): the handback file name returned by the service did not match the
# expected handoff file name. Reflect this in the localization-status report:
#   - Status for that row flips from "Ready for handoff" to
#     "Handback transform failed" (Overview + per-locale sheets).
#   - The per-locale "Error Detail" cell for that row gets the mismatch message.
#   - The "Error Detail" column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn (E) and de-de (F) status columns for the db5902a2 row (row 3)
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale sheets: Status column (C) for the db5902a2 row (row 3)
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Per-locale sheets: Error Detail column (P) for the db5902a2 row (row 3)
$wsZhCn.Range("P3").Value = "Handback file name: yu5ajhxj.mfz is different with handoff file name: db5902a2-193b-479f-b84e-a6b4338caa4f.9c793cdde4a74d5240a718e44983845cf3c1437c.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: yu5ajhxj.mfz is different with handoff file name: db5902a2-193b-479f-b84e-a6b4338caa4f.9c793cdde4a74d5240a718e44983845cf3c1437c.de-de."

# Widen the Error Detail column (P, the 16th column) on both locale sheets to fit
# the new message text (target stored width is 40).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
